# Daten aktualisiert am 2024-01-30
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last used row in column A and append the new tickers after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$newValues = @("IMX-USD", "TAO-USD", "MNT-USD")

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $lastRow + $i + 1
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
